$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    'sliding knee sleeve',
    'volleyball spandex shorts',
    'compression spandex men',
    'volleyball knee pad',
    'knee pads for volleyball',
    'skateboard knee pads',
    'starter youth compression pants',
    'mens compression pants marvel',
    'mouthguard basketball youth',
    'nike compression leggings',
    'knee pad sleeves',
    'jordan flight mens basketball pants',
    'athletic compression pants',
    'nike pro compression leggings men',
    'mens compression tights nike',
    'ladies compression pants',
    'mcdavid compression pants',
    'basketball youth jersey',
    'protec knee pads',
    'skate knee pads',
    'youth knee and elbow pads',
    'under armour compression tights men',
    'elbow and knee pads',
    'adidas tights men',
    'compression pants tesla',
    'poc knee pads',
    'basket knee pads',
    'exercise knee pad',
    'mens compression pants adidas',
    'the rock mens basketball',
    'men''s tights leggings',
    'super compression leggings',
    'under armour compression pants youth boys',
    'knee pad compression pants',
    'pantalones con rodilleras para hombre',
    'basketball tights with knee pads for men',
    'compression pants men basketball',
    'knee leggings for basketball',
    'pants with knee pads for men',
    'knee pad pants for men',
    'mens basketball compression pants',
    'basketball yoga pants',
    'basketball tights with pads',
    'sliding pants mens',
    'basketball compression pants knee pads',
    'compression with knee pads',
    'compression knee pads for men',
    'men basketball compression knee pads',
    'basketball leggings with pads',
    'capri pads',
    'basketball knee pad pants',
    'tights with pads for men',
    'mizuno slider knee pad',
    'mizuno adult slider kneepad',
    'compression knee leggings',
    '5 pad compression shorts',
    'mizuno slider kneepad',
    'youth basketball tights with knee pads',
    'basketball knee pads for men',
    'legging pads',
    'compression capri men pack',
    'softball knee sliding pad',
    'basketball tights knee',
    'compression pads basketball',
    'leggings with knee pads',
    'knee pads baseball',
    'knee length tights men',
    'basketball tight',
    'baseball sliding knee pad',
    'knee compression pants',
    'baseball sliding pad',
    'padded compression pants basketball',
    'compression pants capri men',
    'basketball men leggings',
    'knee sliders softball',
    'mens padded basketball tights',
    'compression pants mens basketball',
    'softball knee slider',
    'mens weightlifting tights',
    'mcdavid 6446 hex knee pads compression leg sleeve',
    'softball sliding pad',
    'pants pad',
    'knee pads mens basketball',
    'padded pants men',
    'basketball legging',
    'padded compression tights basketball',
    'youth tights with knee pads',
    'youth padded tights',
    'baseball knee pads adult',
    'mens capri leggings for sports',
    'basketball tights with knee pads youth boys',
    'youth compression pants with knee pads',
    'mens capris pants',
    'basketball knee pad leggings',
    'padded pants for basketball',
    'compression capri pants men',
    'sliding knee pads',
    'mens capri compression tights',
    'compression below knee',
    'softball sliding pads'
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
